$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-12-18 Thursday"; New = "2025-12-19 Friday" }
    @{ Old = "588×8="; New = "840×8=" }
    @{ Old = "836×6="; New = "630×9=" }
    @{ Old = "499×8="; New = "605×7=" }
    @{ Old = "151×3="; New = "953×7=" }
    @{ Old = "590×5="; New = "386×9=" }
    @{ Old = "401×9="; New = "939×4=" }
    @{ Old = "124×8="; New = "185×9=" }
    @{ Old = "236×5="; New = "542×5=" }
    @{ Old = "311×2="; New = "792×7=" }
    @{ Old = "584×7="; New = "392×7=" }
    @{ Old = "486×3="; New = "675×9=" }
    @{ Old = "479×6="; New = "748×5=" }
    @{ Old = "359×5="; New = "756×3=" }
    @{ Old = "607×6="; New = "115×2=" }
    @{ Old = "817×2="; New = "329×4=" }
    @{ Old = "201×6="; New = "471×8=" }
    @{ Old = "248×6="; New = "949×5=" }
    @{ Old = "802×9="; New = "147×3=" }
    @{ Old = "794×3="; New = "941×9=" }
    @{ Old = "143×6="; New = "427×8=" }
    @{ Old = "252×2="; New = "695×4=" }
    @{ Old = "823×7="; New = "880×6=" }
    @{ Old = "510×5="; New = "758×7=" }
    @{ Old = "733×2="; New = "151×8=" }
    @{ Old = "393×6="; New = "642×7=" }
)

foreach ($pair in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null
}

Write-Output "Replacements applied: $($replacements.Count)"
